# Auto-generated edit script applying the Cerberus_Profits value updates
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets (scheduled market-data refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 50002196
$ws.Range("I62").Value = 50002196
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 50002196
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -50001572
$ws.Range("N62").Value = $null

# Row 65
$ws.Range("H65").Value = 50002196
$ws.Range("I65").Value = 50002196
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 250010980
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -250007860
$ws.Range("N65").Value = $null

# Row 135
$ws.Range("H135").Value = 1560.1428
$ws.Range("I135").Value = 1445.4166
$ws.Range("K135").Value = 13008.7494
$ws.Range("M135").Value = -10473.7494

# Row 138
$ws.Range("H138").Value = 4016.9614
$ws.Range("I138").Value = 4999.4
$ws.Range("J138").Value = 3402.9375
$ws.Range("K138").Value = 14998.2
$ws.Range("L138").Value = 10208.8125
$ws.Range("M138").Value = -9858.199999999999
$ws.Range("N138").Value = -20488.8125

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1095.7273
$ws.Range("I2").Value = 1173.6666
$ws.Range("K2").Value = 1173.6666
$ws.Range("M2").Value = -1060.6666

# Row 61
$ws.Range("H61").Value = 9911
$ws.Range("I61").Value = 6709
$ws.Range("J61").Value = 13913.5
$ws.Range("K61").Value = 6709
$ws.Range("L61").Value = 13913.5
$ws.Range("M61").Value = -6497
$ws.Range("N61").Value = -14337.5

# Row 74
$ws.Range("H74").Value = 5188.5
$ws.Range("I74").Value = 4000.8462
$ws.Range("J74").Value = 6217.8
$ws.Range("K74").Value = 4000.8462
$ws.Range("L74").Value = 6217.8
$ws.Range("M74").Value = -3126.8462
$ws.Range("N74").Value = -7965.8

# Row 77
$ws.Range("H77").Value = 5188.5
$ws.Range("I77").Value = 4000.8462
$ws.Range("J77").Value = 6217.8
$ws.Range("K77").Value = 20004.231
$ws.Range("L77").Value = 31089
$ws.Range("M77").Value = -15636.231
$ws.Range("N77").Value = -39825

# Row 110
$ws.Range("H110").Value = 103444.82
$ws.Range("I110").Value = 125544.22
$ws.Range("K110").Value = 125544.22
$ws.Range("M110").Value = -123499.22

# Row 116
$ws.Range("H116").Value = 1095.7273
$ws.Range("I116").Value = 1173.6666
$ws.Range("K116").Value = 1173.6666
$ws.Range("M116").Value = 1120.3334

# Row 122
$ws.Range("H122").Value = 2844.2
$ws.Range("I122").Value = 2844.2
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8532.599999999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -6082.599999999999
$ws.Range("N122").Value = $null

# Row 134
$ws.Range("H134").Value = 66666.664
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 66666.664
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 66666.664
$ws.Range("M134").Value = $null
$ws.Range("N134").Value = -76806.664

# Row 136
$ws.Range("H136").Value = 9911
$ws.Range("I136").Value = 6709
$ws.Range("J136").Value = 13913.5
$ws.Range("K136").Value = 20127
$ws.Range("L136").Value = 41740.5
$ws.Range("M136").Value = -17577
$ws.Range("N136").Value = -46840.5

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1095.7273
$ws.Range("I3").Value = 1173.6666
$ws.Range("K3").Value = 1173.6666
$ws.Range("M3").Value = -1059.6666

# Row 134
$ws.Range("H134").Value = 11030.15
$ws.Range("I134").Value = 11338.182
$ws.Range("J134").Value = 10653.667
$ws.Range("K134").Value = 34014.546
$ws.Range("L134").Value = 31961.001
$ws.Range("M134").Value = -31479.546
$ws.Range("N134").Value = -37031.001

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 2683.7896
$ws.Range("I58").Value = 1805.5
$ws.Range("J58").Value = 5143
$ws.Range("K58").Value = 1805.5
$ws.Range("L58").Value = 5143
$ws.Range("M58").Value = -1602.5
$ws.Range("N58").Value = -5549

# Row 107
$ws.Range("H107").Value = 3959.0588
$ws.Range("I107").Value = 4246.933
$ws.Range("K107").Value = 4246.933
$ws.Range("M107").Value = -2326.933

# Row 136
$ws.Range("H136").Value = 2683.7896
$ws.Range("I136").Value = 1805.5
$ws.Range("J136").Value = 5143
$ws.Range("K136").Value = 5416.5
$ws.Range("L136").Value = 15429
$ws.Range("M136").Value = -2866.5
$ws.Range("N136").Value = -20529

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 15875469
$ws.Range("I131").Value = 12347291
$ws.Range("J131").Value = 17546712
$ws.Range("K131").Value = 37041873
$ws.Range("L131").Value = 52640136
$ws.Range("M131").Value = -37036833
$ws.Range("N131").Value = -52650216

$ws = $wb.Worksheets.Item("GSM")
# Row 99
$ws.Range("H99").Value = 25899.2
$ws.Range("I99").Value = 8000.5
$ws.Range("J99").Value = 52747.25
$ws.Range("K99").Value = 8000.5
$ws.Range("L99").Value = 52747.25
$ws.Range("M99").Value = -5754.5
$ws.Range("N99").Value = -57239.25

# Row 102
$ws.Range("H102").Value = 4643.5
$ws.Range("J102").Value = 3499
$ws.Range("L102").Value = 3499
$ws.Range("N102").Value = -6743

# Row 107
$ws.Range("H107").Value = 465.63635
$ws.Range("I107").Value = 510.8
$ws.Range("K107").Value = 510.8
$ws.Range("M107").Value = 1409.2

# Row 110
$ws.Range("H110").Value = 35000
$ws.Range("J110").Value = 35000
$ws.Range("L110").Value = 35000
$ws.Range("N110").Value = -43180

# Row 113
$ws.Range("H113").Value = 1455.625
$ws.Range("I113").Value = 1482.8462
$ws.Range("J113").Value = 1337.6666
$ws.Range("K113").Value = 1482.8462
$ws.Range("L113").Value = 1337.6666
$ws.Range("M113").Value = 687.1538
$ws.Range("N113").Value = -5677.6666

# Row 126
$ws.Range("H126").Value = 5710.909
$ws.Range("I126").Value = 2706
$ws.Range("K126").Value = 8118
$ws.Range("M126").Value = -5648

$ws = $wb.Worksheets.Item("LTW")
# Row 74
$ws.Range("H74").Value = 22598
$ws.Range("I74").Value = 22598
$ws.Range("K74").Value = 22598
$ws.Range("M74").Value = -21600

# Row 77
$ws.Range("H77").Value = 22598
$ws.Range("I77").Value = 22598
$ws.Range("K77").Value = 67794
$ws.Range("M77").Value = -62802

# Row 122
$ws.Range("H122").Value = 7733.75
$ws.Range("J122").Value = 7755
$ws.Range("L122").Value = 23265
$ws.Range("N122").Value = -28165

$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 8909.6
$ws.Range("I136").Value = 9565.944
$ws.Range("K136").Value = 28697.832
$ws.Range("M136").Value = -26147.832
